$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the active filter criteria, which unhides all the rows that were
# filtered out (14-18, 21-25, 27-28, 30, 33-37)
$ws.ShowAllData()

# Remove the AutoFilter (and its stored filter/sort state) from the sheet
$ws.AutoFilterMode = $false

# Give column C an explicit width (new col entry in <cols>)
$ws.Range("C1").EntireColumn.ColumnWidth = 10.67

# Move the active selection from H48 to E45
[void]$ws.Range("E45").Select()
